$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new column before column D ("Expected message" / "info" shift right).
$ws.Range("D1").EntireColumn.Insert()
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth()

# New column D header + values ("Type").
$ws.Range("D1").Value = "Type"
$ws.Range("D2").Value = "credentials"
$ws.Range("D3").Value = "credentials"
$ws.Range("D4").Value = "credentials"
$ws.Range("D5").Value = "fieldRequired"
$ws.Range("D6").Value = "fieldRequired"

# Update the (now shifted) "Expected message" column E text.
$ws.Range("E2").Value = "Error: No match for Email and/or Password"
$ws.Range("E3").Value = "Error: No match for Email and/or Password"
$ws.Range("E4").Value = "Error: No match for Email and/or Password"

# Correct the wrong-username values and add matching hyperlinks.
$ws.Range("B3").Value = "wrong@username.com"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:wrong@username.com")
$ws.Range("B3").Style = "Hyperlink"

$ws.Range("B4").Value = "wrong@username.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:wrong@username.com")
$ws.Range("B4").Style = "Hyperlink"

$ws.Activate()
$ws.Range("E13").Select()
